# refactor: rm kafka and worker from CICD
#
# - "Elastic search" row (D8 on the "service" sheet): the curl health-check
#   command moves from https to http.
# - "Kafka connect" row (G5 on the "service" sheet): the
#   kafka-console-consumer command is repointed from the old test topic to
#   the product topic, the stray leading space is dropped, and a
#   key.separator property is appended.
# - "spring cloud" sheet column A gets an explicit custom width (matches the
#   saved file's <cols> entry for column A).

$wb = $excel.ActiveWorkbook

$service = $wb.Worksheets.Item("service")

$service.Range("D8").Value = "curl -k -u elastic:'D_=V-k6LC8zXjpeTPk1V' http://linux-082:50005"
$service.Range("G5").Value = "kafka-console-consumer.sh --bootstrap-server linux-085:50003   --topic mongo.product.products   --from-beginning   --property print.key=true   --property print.value=true   --property key.separator=`" | `""

# Restore the on-screen selection left by the author's last interactive edit.
$service.Activate() | Out-Null
$service.Range("G13").Select() | Out-Null

$springCloud = $wb.Worksheets.Item("spring cloud")
$springCloud.Columns("A:A").ColumnWidth = 13

$springCloud.Activate() | Out-Null
$springCloud.Range("B4").Select() | Out-Null

$service.Activate() | Out-Null
